# Add a /DebugMQTT command and disable MQTT serial print by default.
#
# The target table (Publisher | Subscriber | Topic | Wert | Beschreibung)
# is the 2nd table in the document.  We append three new rows to its end:
#   1. ESP1 -> EM   LastWill/ESP1        "I am going offline"
#   2. EM -> ESP1   ESP1/DebugMQTT  0    "MQTT debugging ausschalten"
#   3. EM -> ESP1   ESP1/DebugMQTT  1    "MQTT debugging einschalten"

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Set-CellRuns($cell, [string[]]$texts) {
    # Write the first chunk directly, then append the remaining chunks
    # just before the cell's end-of-cell marker so the cell ends up
    # containing all chunks concatenated (Word may coalesce adjacent
    # same-formatted text into a single run when the package is saved).
    $cell.Range.Text = $texts[0]
    for ($i = 1; $i -lt $texts.Count; $i++) {
        $full = $cell.Range
        $insertionPoint = $d.Range($full.End - 1, $full.End - 1)
        $insertionPoint.InsertAfter($texts[$i])
    }
}

# ---------------------------------------------------------------------
# Row 1: ESP1 -> EM, topic "LastWill/ESP1"
# ---------------------------------------------------------------------
$row1 = $t.Rows.Add()

Set-CellRuns $row1.Cells.Item(1) @("ESP1")
Set-CellRuns $row1.Cells.Item(2) @("EM")
Set-CellRuns $row1.Cells.Item(3) @("L", "ast", "W", "ill", "/ESP1")
Set-CellRuns $row1.Cells.Item(4) @("I am going offline")
Set-CellRuns $row1.Cells.Item(5) @("ESP1 hat die Verbindung zum MQTT Broaker verloren oder wurde ausgeschlatet.")

# ---------------------------------------------------------------------
# Row 2: EM -> ESP1, topic "ESP1/DebugMQTT", value "0" (disable)
# ---------------------------------------------------------------------
$row2 = $t.Rows.Add()

Set-CellRuns $row2.Cells.Item(1) @("EM")
Set-CellRuns $row2.Cells.Item(2) @("ESP1")
Set-CellRuns $row2.Cells.Item(3) @("ESP1/DebugMQTT")
Set-CellRuns $row2.Cells.Item(4) @("0")
Set-CellRuns $row2.Cells.Item(5) @("MQTT debugging ausschalten")

# ---------------------------------------------------------------------
# Row 3: EM -> ESP1, topic "ESP1/DebugMQTT", value "1" (enable)
# ---------------------------------------------------------------------
$row3 = $t.Rows.Add()

Set-CellRuns $row3.Cells.Item(1) @("EM")
Set-CellRuns $row3.Cells.Item(2) @("ESP1")
Set-CellRuns $row3.Cells.Item(3) @("ESP1/DebugMQTT")
Set-CellRuns $row3.Cells.Item(4) @("1")
Set-CellRuns $row3.Cells.Item(5) @("MQTT debugging", " einschalten")

Write-Output "Table now has $($t.Rows.Count) rows"
